$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 142.4651586029795
$ws.Range("E4").Value = 5.020191085316707
$ws.Range("F4").Value = 0.1991956049093126
$ws.Range("G4").Value = -325.3792547059857
$ws.Range("H4").Value = 11775.03710062947
$ws.Range("I4").Value = 158.8643429998319
$ws.Range("M4").Value = 0.1404877219320042
$ws.Range("N4").Value = 11.85561165795662
$ws.Range("O4").Value = 0.002006714086746797
$ws.Range("P4").Value = 115.829264083106
$ws.Range("Q4").Value = 0.07526798493927345
$ws.Range("R4").Value = 0.351555936911609
$ws.Range("S4").Value = 3253.933034781789
$ws.Range("T4").Value = -117738.5153946368
$ws.Range("U4").Value = -1588.641423284233
$ws.Range("V4").Value = -280.8629212975648
$ws.Range("W4").Value = -2704.695545919531
$ws.Range("X4").Value = -3777.110015785904
